$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Ativação:" date string 01/01/2020 -> 01/01/2022 (kept as literal TEXT, not an
#    Excel date serial). Typing a plain date-shaped string into B8/C8 would make
#    Excel auto-convert it to a date value + mint a new number-format style, which
#    does not match the source (string stays plain text, same cell style). Route the
#    literal through a text formula and paste back as a value to avoid the date
#    auto-detect while keeping the original style untouched.
$ws.Range("B8:C8").Formula = "=""01/01/2022"""
$ws.Range("B8:C8").Copy()
$ws.Range("B8:C8").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# 2) "Programa:" row (16) gets the full, detailed Portuguese programme text
#    (previously it duplicated the short "Programa resumido" text).
$programaDetalhado = "- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores."
$ws.Range("B16:C16").Value = $programaDetalhado

# 3) "Syllabus:" row (17) gets the full, detailed English syllabus text
#    (previously it duplicated the "Short syllabus" text).
$syllabusDetalhado = "- fluid rheology,- Sizing of pipes,- Accessories and pumping for industrial fluids,- Stirring and mixing,- Particle characterization and particle bed,- Sedimentation,- Filtration,- Processes with membranes.- Unit heat exchange operations: heat exchangers and evaporators."
$ws.Range("B17:C17").Value = $syllabusDetalhado
